$d = $word.ActiveDocument
$d.Content.Find.Execute("2022-12-24 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2022-12-25 Sunday", 2) | Out-Null
$d.Content.Find.Execute("47-28=", $true, $false, $false, $false, $false, $true, 1, $false, "0+93=", 2) | Out-Null
$d.Content.Find.Execute("43+39=", $true, $false, $false, $false, $false, $true, 1, $false, "73+11=", 2) | Out-Null
$d.Content.Find.Execute("0+26=", $true, $false, $false, $false, $false, $true, 1, $false, "80-51=", 2) | Out-Null
$d.Content.Find.Execute("34+40=", $true, $false, $false, $false, $false, $true, 1, $false, "39-5=", 2) | Out-Null
$d.Content.Find.Execute("0+88=", $true, $false, $false, $false, $false, $true, 1, $false, "51-33=", 2) | Out-Null
$d.Content.Find.Execute("26-19=", $true, $false, $false, $false, $false, $true, 1, $false, "29+61=", 2) | Out-Null
$d.Content.Find.Execute("82-78=", $true, $false, $false, $false, $false, $true, 1, $false, "91-62=", 2) | Out-Null
$d.Content.Find.Execute("82-79=", $true, $false, $false, $false, $false, $true, 1, $false, "57-8=", 2) | Out-Null
$d.Content.Find.Execute("50+23=", $true, $false, $false, $false, $false, $true, 1, $false, "56+13=", 2) | Out-Null
$d.Content.Find.Execute("43+18=", $true, $false, $false, $false, $false, $true, 1, $false, "89-10=", 2) | Out-Null
$d.Content.Find.Execute("59+4=", $true, $false, $false, $false, $false, $true, 1, $false, "70-16=", 2) | Out-Null
$d.Content.Find.Execute("17+70=", $true, $false, $false, $false, $false, $true, 1, $false, "97-4=", 2) | Out-Null
$d.Content.Find.Execute("46+52=", $true, $false, $false, $false, $false, $true, 1, $false, "10+11=", 2) | Out-Null
$d.Content.Find.Execute("23+45=", $true, $false, $false, $false, $false, $true, 1, $false, "70-16=", 2) | Out-Null
$d.Content.Find.Execute("79+20=", $true, $false, $false, $false, $false, $true, 1, $false, "25+6=", 2) | Out-Null
$d.Content.Find.Execute("7+53=", $true, $false, $false, $false, $false, $true, 1, $false, "83-49=", 2) | Out-Null
$d.Content.Find.Execute("47-16=", $true, $false, $false, $false, $false, $true, 1, $false, "78-4=", 2) | Out-Null
$d.Content.Find.Execute("29+49=", $true, $false, $false, $false, $false, $true, 1, $false, "12+74=", 2) | Out-Null
$d.Content.Find.Execute("2+69=", $true, $false, $false, $false, $false, $true, 1, $false, "65-58=", 2) | Out-Null
$d.Content.Find.Execute("71-61=", $true, $false, $false, $false, $false, $true, 1, $false, "33-5=", 2) | Out-Null
$d.Content.Find.Execute("46+18=", $true, $false, $false, $false, $false, $true, 1, $false, "35+52=", 2) | Out-Null
$d.Content.Find.Execute("16+79=", $true, $false, $false, $false, $false, $true, 1, $false, "18+32=", 2) | Out-Null
$d.Content.Find.Execute("95-23=", $true, $false, $false, $false, $false, $true, 1, $false, "90-69=", 2) | Out-Null
$d.Content.Find.Execute("19+36=", $true, $false, $false, $false, $false, $true, 1, $false, "85-78=", 2) | Out-Null
$d.Content.Find.Execute("57+38=", $true, $false, $false, $false, $false, $true, 1, $false, "98-61=", 2) | Out-Null
$d.Content.Find.Execute("7+89=", $true, $false, $false, $false, $false, $true, 1, $false, "61-16=", 2) | Out-Null
$d.Content.Find.Execute("5+21=", $true, $false, $false, $false, $false, $true, 1, $false, "70+2=", 2) | Out-Null
$d.Content.Find.Execute("71+26=", $true, $false, $false, $false, $false, $true, 1, $false, "36+37=", 2) | Out-Null
$d.Content.Find.Execute("36+34=", $true, $false, $false, $false, $false, $true, 1, $false, "21+9=", 2) | Out-Null
$d.Content.Find.Execute("21+70=", $true, $false, $false, $false, $false, $true, 1, $false, "67+30=", 2) | Out-Null
$d.Content.Find.Execute("78-2=", $true, $false, $false, $false, $false, $true, 1, $false, "65+25=", 2) | Out-Null
$d.Content.Find.Execute("57+36=", $true, $false, $false, $false, $false, $true, 1, $false, "72-71=", 2) | Out-Null
$d.Content.Find.Execute("25-18=", $true, $false, $false, $false, $false, $true, 1, $false, "26+34=", 2) | Out-Null
$d.Content.Find.Execute("84-76=", $true, $false, $false, $false, $false, $true, 1, $false, "50+33=", 2) | Out-Null
$d.Content.Find.Execute("45+40=", $true, $false, $false, $false, $false, $true, 1, $false, "78-72=", 2) | Out-Null
$d.Content.Find.Execute("22-18=", $true, $false, $false, $false, $false, $true, 1, $false, "23-9=", 2) | Out-Null
$d.Content.Find.Execute("56+19=", $true, $false, $false, $false, $false, $true, 1, $false, "77+18=", 2) | Out-Null
$d.Content.Find.Execute("97-40=", $true, $false, $false, $false, $false, $true, 1, $false, "8+8=", 2) | Out-Null
$d.Content.Find.Execute("23+15=", $true, $false, $false, $false, $false, $true, 1, $false, "65-24=", 2) | Out-Null
$d.Content.Find.Execute("30-7=", $true, $false, $false, $false, $false, $true, 1, $false, "22+13=", 2) | Out-Null
$d.Content.Find.Execute("83-28=", $true, $false, $false, $false, $false, $true, 1, $false, "0+75=", 2) | Out-Null
$d.Content.Find.Execute("54-36=", $true, $false, $false, $false, $false, $true, 1, $false, "62-43=", 2) | Out-Null
$d.Content.Find.Execute("86-4=", $true, $false, $false, $false, $false, $true, 1, $false, "77-49=", 2) | Out-Null
$d.Content.Find.Execute("10+71=", $true, $false, $false, $false, $false, $true, 1, $false, "28+64=", 2) | Out-Null
$d.Content.Find.Execute("60-55=", $true, $false, $false, $false, $false, $true, 1, $false, "95-39=", 2) | Out-Null
$d.Content.Find.Execute("53+34=", $true, $false, $false, $false, $false, $true, 1, $false, "78-48=", 2) | Out-Null
$d.Content.Find.Execute("73+19=", $true, $false, $false, $false, $false, $true, 1, $false, "44+20=", 2) | Out-Null
$d.Content.Find.Execute("88-58=", $true, $false, $false, $false, $false, $true, 1, $false, "9+3=", 2) | Out-Null
$d.Content.Find.Execute("33+60=", $true, $false, $false, $false, $false, $true, 1, $false, "93-6=", 2) | Out-Null
$d.Content.Find.Execute("55-16=", $true, $false, $false, $false, $false, $true, 1, $false, "65-5=", 2) | Out-Null
$d.Content.Find.Execute("90-10=", $true, $false, $false, $false, $false, $true, 1, $false, "21+33=", 2) | Out-Null
$d.Content.Find.Execute("0+34=", $true, $false, $false, $false, $false, $true, 1, $false, "75-55=", 2) | Out-Null
$d.Content.Find.Execute("21+26=", $true, $false, $false, $false, $false, $true, 1, $false, "68-55=", 2) | Out-Null
$d.Content.Find.Execute("27+5=", $true, $false, $false, $false, $false, $true, 1, $false, "71-24=", 2) | Out-Null
$d.Content.Find.Execute("24+18=", $true, $false, $false, $false, $false, $true, 1, $false, "68-62=", 2) | Out-Null
$d.Content.Find.Execute("49-11=", $true, $false, $false, $false, $false, $true, 1, $false, "93-90=", 2) | Out-Null
$d.Content.Find.Execute("44-10=", $true, $false, $false, $false, $false, $true, 1, $false, "15+71=", 2) | Out-Null
$d.Content.Find.Execute("68-41=", $true, $false, $false, $false, $false, $true, 1, $false, "35+46=", 2) | Out-Null
$d.Content.Find.Execute("64-13=", $true, $false, $false, $false, $false, $true, 1, $false, "34+49=", 2) | Out-Null
$d.Content.Find.Execute("34+11=", $true, $false, $false, $false, $false, $true, 1, $false, "3+92=", 2) | Out-Null
$d.Content.Find.Execute("28-9=", $true, $false, $false, $false, $false, $true, 1, $false, "10+10=", 2) | Out-Null
$d.Content.Find.Execute("89-25=", $true, $false, $false, $false, $false, $true, 1, $false, "78-42=", 2) | Out-Null
$d.Content.Find.Execute("19+38=", $true, $false, $false, $false, $false, $true, 1, $false, "33+2=", 2) | Out-Null
$d.Content.Find.Execute("4+58=", $true, $false, $false, $false, $false, $true, 1, $false, "81-27=", 2) | Out-Null
$d.Content.Find.Execute("50-5=", $true, $false, $false, $false, $false, $true, 1, $false, "49-48=", 2) | Out-Null
$d.Content.Find.Execute("56+32=", $true, $false, $false, $false, $false, $true, 1, $false, "50-31=", 2) | Out-Null
$d.Content.Find.Execute("84-23=", $true, $false, $false, $false, $false, $true, 1, $false, "17-5=", 2) | Out-Null
$d.Content.Find.Execute("47+9=", $true, $false, $false, $false, $false, $true, 1, $false, "52+47=", 2) | Out-Null
$d.Content.Find.Execute("75-23=", $true, $false, $false, $false, $false, $true, 1, $false, "87-25=", 2) | Out-Null
$d.Content.Find.Execute("67-41=", $true, $false, $false, $false, $false, $true, 1, $false, "18+8=", 2) | Out-Null
$d.Content.Find.Execute("27+49=", $true, $false, $false, $false, $false, $true, 1, $false, "54+9=", 2) | Out-Null
$d.Content.Find.Execute("49-33=", $true, $false, $false, $false, $false, $true, 1, $false, "57+40=", 2) | Out-Null
$d.Content.Find.Execute("0+73=", $true, $false, $false, $false, $false, $true, 1, $false, "55+28=", 2) | Out-Null
$d.Content.Find.Execute("55+23=", $true, $false, $false, $false, $false, $true, 1, $false, "46+43=", 2) | Out-Null
$d.Content.Find.Execute("11+22=", $true, $false, $false, $false, $false, $true, 1, $false, "50-15=", 2) | Out-Null
$d.Content.Find.Execute("39+39=", $true, $false, $false, $false, $false, $true, 1, $false, "90-12=", 2) | Out-Null
$d.Content.Find.Execute("36-11=", $true, $false, $false, $false, $false, $true, 1, $false, "4+41=", 2) | Out-Null
$d.Content.Find.Execute("65+24=", $true, $false, $false, $false, $false, $true, 1, $false, "49+1=", 2) | Out-Null
$d.Content.Find.Execute("98-47=", $true, $false, $false, $false, $false, $true, 1, $false, "85-16=", 2) | Out-Null
$d.Content.Find.Execute("42+41=", $true, $false, $false, $false, $false, $true, 1, $false, "63+11=", 2) | Out-Null
$d.Content.Find.Execute("85-8=", $true, $false, $false, $false, $false, $true, 1, $false, "45-2=", 2) | Out-Null
$d.Content.Find.Execute("86-32=", $true, $false, $false, $false, $false, $true, 1, $false, "93-56=", 2) | Out-Null
$d.Content.Find.Execute("8+6=", $true, $false, $false, $false, $false, $true, 1, $false, "95-77=", 2) | Out-Null
$d.Content.Find.Execute("68-15=", $true, $false, $false, $false, $false, $true, 1, $false, "8+29=", 2) | Out-Null
$d.Content.Find.Execute("30+67=", $true, $false, $false, $false, $false, $true, 1, $false, "21-8=", 2) | Out-Null
$d.Content.Find.Execute("36+10=", $true, $false, $false, $false, $false, $true, 1, $false, "49+25=", 2) | Out-Null
$d.Content.Find.Execute("24+66=", $true, $false, $false, $false, $false, $true, 1, $false, "86-63=", 2) | Out-Null
$d.Content.Find.Execute("47+1=", $true, $false, $false, $false, $false, $true, 1, $false, "46-29=", 2) | Out-Null
$d.Content.Find.Execute("15+0=", $true, $false, $false, $false, $false, $true, 1, $false, "2+75=", 2) | Out-Null
$d.Content.Find.Execute("41+0=", $true, $false, $false, $false, $false, $true, 1, $false, "6+84=", 2) | Out-Null
$d.Content.Find.Execute("49+23=", $true, $false, $false, $false, $false, $true, 1, $false, "16+69=", 2) | Out-Null
$d.Content.Find.Execute("39+32=", $true, $false, $false, $false, $false, $true, 1, $false, "29+43=", 2) | Out-Null
$d.Content.Find.Execute("47+41=", $true, $false, $false, $false, $false, $true, 1, $false, "79+13=", 2) | Out-Null
$d.Content.Find.Execute("53+3=", $true, $false, $false, $false, $false, $true, 1, $false, "56+11=", 2) | Out-Null
$d.Content.Find.Execute("31+40=", $true, $false, $false, $false, $false, $true, 1, $false, "84-71=", 2) | Out-Null
$d.Content.Find.Execute("58-9=", $true, $false, $false, $false, $false, $true, 1, $false, "94-58=", 2) | Out-Null
$d.Content.Find.Execute("20-12=", $true, $false, $false, $false, $false, $true, 1, $false, "17+30=", 2) | Out-Null
$d.Content.Find.Execute("59-13=", $true, $false, $false, $false, $false, $true, 1, $false, "49+2=", 2) | Out-Null
$d.Content.Find.Execute("99-86=", $true, $false, $false, $false, $false, $true, 1, $false, "91-28=", 2) | Out-Null
$d.Content.Find.Execute("59-6=", $true, $false, $false, $false, $false, $true, 1, $false, "67-39=", 2) | Out-Null
